$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 3 (new weekly price observation), pushing
# the existing rows 3-9 down to rows 4-10.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44881
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103003
$ws.Range("J3").Value = "Damasco"
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 41000
$ws.Range("O3").Value = 42000
$ws.Range("P3").Value = 41500
$ws.Range("Q3").Value = "$/bandeja 18 kilos"
$ws.Range("R3").Value = "Región de Coquimbo"
$ws.Range("S3").Value = 2306
$ws.Range("T3").Value = 18
